$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns whose values need to be swapped between row 5 and row 6.
# (Other columns on these rows contain identical values in both rows,
# so no change is needed there.)
$cols = @("A","B","D","E","F","G","H","I","Q","R")

foreach ($col in $cols) {
    $addr5 = $col + "5"
    $addr6 = $col + "6"

    $cell5 = $ws.Range($addr5)
    $cell6 = $ws.Range($addr6)

    # Value2 returns the raw (un-formatted) value and reads correctly
    # into a variable in this runtime.
    $v5 = $cell5.Value2
    $v6 = $cell6.Value2

    if ($col -eq "I") {
        # Column I stores numeric-looking text ("5" / "1") as text, not
        # numbers. Force text formatting before assigning so Excel
        # doesn't auto-convert the string back into a number, then
        # restore the default style so no extra formatting is left
        # behind.
        $cell5.NumberFormat = "@"
        $cell5.Value = [string]$v6
        $cell5.Style = "Normal"

        $cell6.NumberFormat = "@"
        $cell6.Value = [string]$v5
        $cell6.Style = "Normal"
    } else {
        $cell5.Value = $v6
        $cell6.Value = $v5
    }
}
